# Applies the rcbus-opl3 BoM update:
#  - C2 joins the 10uF polarized-cap group (C8 C15 C16 -> C2 C8 C15 C16)
#  - That cap's value text is re-cased (10uF -> 10uf)
#  - Quantity Per PCB for that group goes from 3 to 4
#  - The board oscillator part number is replaced by its frequency (SG-8002CA -> 14.3181MHz)
#  - Component/Total counts bump from 31 to 32 (29 SMD/ 3 THT)
#  - Regenerated-on date / timestamp move forward

$wb = $excel.ActiveWorkbook

$bom = $wb.Worksheets.Item("BoM")
$costs = $wb.Worksheets.Item("Costs")

# --- BoM sheet: polarized capacitor group row (row 12) ---
$bom.Range("D12").Value = "C2 C8 C15 C16"
$bom.Range("E12").Value = "10uf"
$bom.Range("G12").Value = "4"

# --- BoM sheet: oscillator row (row 24) value column ---
$bom.Range("E24").Value = "14.3181MHz"

# --- BoM sheet: header summary counts ---
$bom.Range("F6").Value = 32

# --- Costs sheet: matching quantity formula for the same cap group (row 13) ---
$costs.Range("G13").Formula = "=CEILING(BoardQty*4,1)"

# --- Costs sheet: component-count summary string ---
$costs.Range("E6").Value = "32 (29 SMD/ 3 THT)"

# --- Costs sheet: regeneration date/time ---
$costs.Range("E5").Value = "2025-10-04"
$costs.Range("B28").Value = "2025-10-04 22:33:57"
